# Working on decision trees: relabel the "target" rows to "met"/"not met"
# on both the Perf and Class sheets, and leave the selection on E9 with
# the Perf sheet as the active tab.

$wb = $excel.ActiveWorkbook

$wsPerf  = $wb.Worksheets.Item("Perf")
$wsClass = $wb.Worksheets.Item("Class")

foreach ($ws in @($wsPerf, $wsClass)) {
    $ws.Range("E4").Value = "Left met"
    $ws.Range("E5").Value = "Left not met"
    $ws.Range("E7").Value = "Right met"
    $ws.Range("E8").Value = "Right not met"
}

# Leave the Class sheet selected on E9 first ...
[void]$wsClass.Select()
[void]$wsClass.Range("E9").Select()

# ... then switch to (and leave active on) the Perf sheet, also at E9.
[void]$wsPerf.Select()
[void]$wsPerf.Range("E9").Select()
